$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (12 and 13) with the same data pattern as the existing
# rows (A: 102, B: nokia, C: 999), for the newly added UT and FT test
# cases (authentication / integration).
$ws.Range("A2:C2").Copy()
$ws.Range("A12:C12").PasteSpecial()

$ws.Range("A2:C2").Copy()
$ws.Range("A13:C13").PasteSpecial()
